$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '43.840.24'
$ws.Range("E2").Value = '  +1.66%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.284.23'
$ws.Range("E3").Value = '  +1.24%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.53%  '

# Row 5: Solana
$ws.Range("D5").Value = '''120.78'
$ws.Range("E5").Value = '  +7.49%  '

# Row 6: BNB
$ws.Range("D6").Value = '''266.90'
$ws.Range("E6").Value = '  +1.25%  '

# Row 7: XRP
$ws.Range("E7").Value = '  +4.82%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.07%  '

# Row 9: Cardano
$ws.Range("D9").Value = '''0.628'
$ws.Range("E9").Value = '  +4.92%  '

# Row 10: Avalanche
$ws.Range("D10").Value = '''48.43'
$ws.Range("E10").Value = '  +2.38%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = '''0.0952'
$ws.Range("E11").Value = '  +3.50%  '

# Row 12: Polkadot
$ws.Range("D12").Value = '''9.27'
$ws.Range("E12").Value = '  +7.03%  '

# Row 13: TRON
$ws.Range("D13").Value = '''0.107'
$ws.Range("E13").Value = '  -0.16%  '

# Row 14: Chainlink
$ws.Range("D14").Value = '''15.63'
$ws.Range("E14").Value = '  +1.87%  '

# Row 15: Polygon
$ws.Range("D15").Value = '''0.920'
$ws.Range("E15").Value = '  +8.87%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = '2.620.90'
$ws.Range("E16").Value = '  +0.16%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '2.285.70'
$ws.Range("E17").Value = '  +0.69%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '43.761.06'
$ws.Range("E18").Value = '  +1.17%  '

# Row 19: ShibaInu
$ws.Range("D19").Value = '''0.0000111'
$ws.Range("E19").Value = '  +3.94%  '

# Row 20: Uniswap
$ws.Range("D20").Value = '''6.95'
$ws.Range("E20").Value = '  +1.41%  '

# Row 21: Litecoin
$ws.Range("D21").Value = '''72.41'
$ws.Range("E21").Value = '  +2.36%  '

# Row 22: ImmutableX
$ws.Range("D22").Value = '''2.42'
$ws.Range("E22").Value = '  +2.38%  '

# Row 23: BitcoinCash
$ws.Range("D23").Value = '''236.86'
$ws.Range("E23").Value = '  +3.61%  '

# Row 24: InternetComputer(DFINITY)
$ws.Range("D24").Value = '''9.67'
$ws.Range("E24").Value = '  +0.04%  '

# Row 25: PancakeSwap
$ws.Range("D25").Value = '''2.89'
$ws.Range("E25").Value = '  +2.60%  '

# Row 26: Cosmos
$ws.Range("D26").Value = '''12.05'
$ws.Range("E26").Value = '  +7.22%  '

# Row 27: Dai
$ws.Range("E27").Value = '  +1.92%  '

# Row 28: InjectiveProtocol
$ws.Range("D28").Value = '''42.47'
$ws.Range("E28").Value = '  +4.77%  '

# Row 29: WEMIXToken
$ws.Range("D29").Value = '''3.38'
$ws.Range("E29").Value = '  -0.28%  '

# Row 30: Toncoin
$ws.Range("E30").Value = '  +0.30%  '

# Row 31: Monero
$ws.Range("D31").Value = '''172.59'
$ws.Range("E31").Value = '  +0.62%  '

# Row 32: Hedera
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.0934'
$ws.Range("E32").Value = '  +4.31%  '

# Row 33: EthereumClassic
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''21.70'
$ws.Range("E33").Value = '  +2.97%  '

# Row 34: Filecoin
$ws.Range("D34").Value = '''5.79'
$ws.Range("E34").Value = '  +5.21%  '

# Row 35: Stellar
$ws.Range("E35").Value = '  +4.44%  '

# Row 36: NEARProtocol
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '''4.27'
$ws.Range("E36").Value = '  +15.46%  '

# Row 37: VeChain
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.0385'
$ws.Range("E37").Value = '  +11.77%  '

# Row 38: RenderToken
$ws.Range("D38").Value = '''4.63'
$ws.Range("E38").Value = '  +1.83%  '

# Row 39: Kaspa
$ws.Range("D39").Value = '''0.107'
$ws.Range("E39").Value = '  +3.80%  '

# Row 40: LidoDAOToken
$ws.Range("D40").Value = '''2.58'
$ws.Range("E40").Value = '  +8.56%  '

# Row 41: MultiversX
$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").Value = '''74.42'
$ws.Range("E41").Value = '  +1.54%  '

# Row 42: Celestia
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").Value = '''13.91'
$ws.Range("E42").Value = '  -1.24%  '

# Row 43: Algorand
$ws.Range("D43").Value = '''0.239'
$ws.Range("E43").Value = '  +3.13%  '

# Row 44: FirstDigitalUSD
$ws.Range("E44").Value = '  -0.46%  '

# Row 45: ARBITRUM
$ws.Range("D45").Value = '''1.39'
$ws.Range("E45").Value = '  +2.86%  '

# Row 46: THORChain
$ws.Range("D46").Value = '''5.84'
$ws.Range("E46").Value = '  -4.07%  '

# Row 47: ordi
$ws.Range("D47").Value = '''76.80'
$ws.Range("E47").Value = '  +51.66%  '

# Row 48: TrustWalletToken
$ws.Range("E48").Value = '  +3.92%  '

# Row 49: FraxShare
$ws.Range("D49").Value = '''8.58'
$ws.Range("E49").Value = '  +0.58%  '

# Row 50: Cronos
$ws.Range("E50").Value = '  +2.34%  '

# Row 51: Aave
$ws.Range("D51").Value = '''102.52'
$ws.Range("E51").Value = '  +3.21%  '
